# Switches-BOM.xlsx update:
#  - H1 H5 hardware group gains H6 and H10 (qty 2 -> 4)
#  - H3 hardware group gains H8 (qty 1 -> 2)
#  - cursor/selection left on C7 (the updated standoff row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Screws row (was "H1 H5", qty 2/2)
$ws.Range("A6").Value = "H1 H5 H6 H10"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 4

# Standoff row (was "H3", qty 1/1)
$ws.Range("A7").Value = "H3 H8"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2

# Leave the active selection on the standoff row, matching the saved file
$ws.Range("C7").Select()
